$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.25
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.75
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 9
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 41
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 21
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 2.63
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 26
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 126
$ws.Range("BC2").Value = 151

$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 2.63
$ws.Range("L3").Value = 6.5
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.2
$ws.Range("X3").Value = 6.5
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 21
$ws.Range("AK3").Value = 67
$ws.Range("AM3").Value = 81
$ws.Range("AN3").Value = 3.5

$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5

$ws.Range("H6").Value = 2.88
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 1.83
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 2.25
$ws.Range("Q6").Value = 2.88
$ws.Range("R6").Value = 1.4
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 6
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 26
$ws.Range("AF6").Value = 81
$ws.Range("AH6").Value = 7
$ws.Range("AU6").Value = 9.5
$ws.Range("AZ6").Value = 67

$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 4.2
$ws.Range("I7").Value = 6.25
$ws.Range("J7").Value = 2.05
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("Q7").Value = 1.83
$ws.Range("R7").Value = 2.03
$ws.Range("Z7").Value = 10
$ws.Range("AD7").Value = 8
$ws.Range("AE7").Value = 19
$ws.Range("AN7").Value = 3.4

$ws.Range("G8").Value = 3.2
$ws.Range("I8").Value = 2.45
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 1.95
$ws.Range("L8").Value = 3.25
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.63
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 1.53
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("W8").Value = 8
$ws.Range("Y8").Value = 12
$ws.Range("AD8").Value = 6
$ws.Range("AE8").Value = 17
$ws.Range("AI8").Value = 11
$ws.Range("AK8").Value = 23
$ws.Range("AN8").Value = 4.75
$ws.Range("AX8").Value = 15
$ws.Range("AY8").Value = 29

$ws.Range("I9").Value = 3.7
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("AH9").Value = 10
